# Adding trend zone violations: replace a block of noisy raw values with
# values that push the series past the trend-rule control limits, and
# reset the sheet's view/selection back to the top of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value  = -3
$ws.Range("B10").Value = -2
$ws.Range("B11").Value = 0.2
$ws.Range("B12").Value = 0.7
$ws.Range("B14").Value = 1.2

# Reset the view: scroll back to the top and select B15.
$ws.Range("B15").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
